# [Fonds de solidarite] Add 2020-08-17 data
# Updates nombre_aides (col C) and montant_total (col D) for the VOLET2
# regional / classe_effectif rows whose underlying counts changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; C = "178"; D = "408816.00" },
    @{ Row = 3; C = "954"; D = "2949159.33" },
    @{ Row = 4; C = "401"; D = "1611011.12" },
    @{ Row = 5; C = "106"; D = "476466.09" },
    @{ Row = 7; C = "6"; D = "35500.00" },
    @{ Row = 8; C = "33"; D = "66000.00" },
    @{ Row = 9; C = "55"; D = "135928.41" },
    @{ Row = 10; C = "341"; D = "1143495.69" },
    @{ Row = 11; C = "140"; D = "576601.77" },
    @{ Row = 12; C = "34"; D = "165120.00" },
    @{ Row = 13; C = "7"; D = "37000.00" },
    @{ Row = 15; C = "100"; D = "271752.38" },
    @{ Row = 16; C = "424"; D = "1329636.54" },
    @{ Row = 17; C = "150"; D = "635869.10" },
    @{ Row = 19; C = "17"; D = "111216.00" },
    @{ Row = 33; C = "104"; D = "283673.00" },
    @{ Row = 34; C = "551"; D = "1760650.66" },
    @{ Row = 35; C = "220"; D = "1088288.11" },
    @{ Row = 36; C = "72"; D = "387894.00" },
    @{ Row = 39; C = "34"; D = "83830.00" },
    @{ Row = 40; C = "154"; D = "391828.00" },
    @{ Row = 41; C = "81"; D = "276298.00" },
    @{ Row = 42; C = "20"; D = "77995.14" },
    @{ Row = 44; C = "49"; D = "110683.00" },
    @{ Row = 50; C = "97"; D = "273768.17" },
    @{ Row = 51; C = "560"; D = "1904876.52" },
    @{ Row = 52; C = "256"; D = "1096446.76" },
    @{ Row = 55; C = "18"; D = "56720.65" },
    @{ Row = 56; C = "709"; D = "1807318.62" },
    @{ Row = 57; C = "3503"; D = "10921054.21" },
    @{ Row = 58; C = "1814"; D = "7529834.94" },
    @{ Row = 59; C = "627"; D = "3066791.45" },
    @{ Row = 60; C = "127"; D = "875123.00" },
    @{ Row = 62; C = "284"; D = "674120.58" },
    @{ Row = 80; C = "875"; D = "2767444.67" },
    @{ Row = 81; C = "332"; D = "1321403.79" },
    @{ Row = 82; C = "116"; D = "561984.52" },
    @{ Row = 97; C = "281"; D = "725779.43" },
    @{ Row = 98; C = "1183"; D = "3561000.70" },
    @{ Row = 99; C = "442"; D = "1762594.02" },
    @{ Row = 101; C = "31"; D = "189157.00" }
)

foreach ($u in $updates) {
    $row = $u.Row
    # Values are text (inline strings) in the source data, e.g. "408816.00",
    # so force Text format before assigning to keep them from being
    # reinterpreted as numbers (which would drop the trailing zeros).
    $ws.Range("C$row" + ":D$row").NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $u.C
    $ws.Cells.Item($row, 4).Value = $u.D
}
